$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.62175190448761
$ws.Range("B1").Value = -1
$ws.Range("C1").Value = 2.647257566452026
$ws.Range("D1").Value = 1.395869731903076
$ws.Range("E1").Value = 0.9958027601242065
